$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.6353409687593692
$ws.Range("C2").Value = -0.01857434322460993
$ws.Range("D2").Value = 1.312052523820988
$ws.Range("E2").Value = 454
$ws.Range("F2").Value = 729

$ws.Range("B3").Value = -0.195308571777268
$ws.Range("C3").Value = 0.001434192486353014
$ws.Range("D3").Value = 1.182590850842198
$ws.Range("E3").Value = 573
$ws.Range("F3").Value = 731
$ws.Range("G3").Value = 1044

$ws.Range("B4").Value = -0.05519750408199842
$ws.Range("C4").Value = -0.003339941030181426
$ws.Range("D4").Value = 0.04439901661719293
$ws.Range("E4").Value = 727
$ws.Range("F4").Value = 577

$ws.Range("B5").Value = 0.1628605380343664
$ws.Range("C5").Value = -0.007243684334445333
$ws.Range("D5").Value = 0.02843900994477311
$ws.Range("E5").Value = 571
$ws.Range("F5").Value = 733

$ws.Range("B6").Value = 0.5026690931498636
$ws.Range("C6").Value = -0.0009943161996567307
$ws.Range("D6").Value = 0.005409761582878377
$ws.Range("E6").Value = 754
$ws.Range("F6").Value = 550

$ws.Range("B7").Value = -0.1027401287353118
$ws.Range("C7").Value = -0.0006036375541214056
$ws.Range("D7").Value = 0.1262990093548035
$ws.Range("E7").Value = 612
$ws.Range("F7").Value = 692

$ws.Range("B8").Value = -0.1959515989499447
$ws.Range("C8").Value = 0.001394404822410822
$ws.Range("D8").Value = 0.6806948648731734
$ws.Range("E8").Value = 669
$ws.Range("F8").Value = 635

$ws.Range("B9").Value = 0.06596824732900508
$ws.Range("C9").Value = 0.0006162046821136524
$ws.Range("D9").Value = 0.6960589546184718
$ws.Range("E9").Value = 629
$ws.Range("F9").Value = 675
$ws.Range("G9").Value = 1548

$ws.Range("B10").Value = -0.1982798011535828
$ws.Range("C10").Value = 0.001036835673354775
$ws.Range("D10").Value = 1.30595720970737
$ws.Range("E10").Value = 568
$ws.Range("F10").Value = 736

$ws.Range("B11").Value = -0.1162315210713972
$ws.Range("C11").Value = -0.001335489015084623
$ws.Range("D11").Value = 0.1142483992538074
$ws.Range("E11").Value = 525
$ws.Range("F11").Value = 501

$ws.Range("B12").Value = -0.2248178812031053
$ws.Range("C12").Value = 0.01503158482983724
$ws.Range("D12").Value = 0.3651128817038775
$ws.Range("E12").Value = 598
$ws.Range("F12").Value = 706
$ws.Range("G12").Value = 1606

$ws.Range("B13").Value = -0.1544684014222742
$ws.Range("C13").Value = 0.0006416574988037382
$ws.Range("D13").Value = 0.02146678169534579
$ws.Range("E13").Value = 662
$ws.Range("F13").Value = 642
$ws.Range("H13").Value = 1343

$ws.Range("B14").Value = -0.08078688171742576
$ws.Range("C14").Value = -0.001343924020154708
$ws.Range("D14").Value = 0.9267364681758432
$ws.Range("E14").Value = 628
$ws.Range("F14").Value = 676

$ws.Range("B15").Value = -0.2330844633052901
$ws.Range("C15").Value = -0.001281915713442383
$ws.Range("D15").Value = 0.3059710752177072
$ws.Range("E15").Value = 524
$ws.Range("F15").Value = 780

$ws.Range("B16").Value = -0.1654167502093068
$ws.Range("C16").Value = 0.0006593113420236829
$ws.Range("D16").Value = 0.2437648469626392
$ws.Range("E16").Value = 554
$ws.Range("F16").Value = 631

$ws.Range("B17").Value = 0.4919399178650536
$ws.Range("C17").Value = 0.0120796188593173
$ws.Range("D17").Value = 0.1235119339359379
$ws.Range("E17").Value = 523
$ws.Range("F17").Value = 781
